# mejore el uso de directorios
# Inserts 7 new container-withdrawal records into the "retiros_puerto" sheet,
# shifting the existing 57 data rows down to make room. After the inserts,
# the sheet grows from A1:E58 to A1:E65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel row numbers (in the FINAL sheet) at which a brand-new row must appear.
# Inserting at these positions in ascending order reproduces the target
# layout exactly, since each insertion only pushes rows at-or-below it down
# by one and never disturbs earlier, already-placed rows.
$insertRows = @(2, 3, 4, 5, 6, 9, 40)

foreach ($r in $insertRows) {
    $ws.Rows.Item($r).Insert()
    # Insert() copies the formatting of the row above (here, sometimes the
    # bold/bordered header row), which the source data rows never had.
    # Reset to the workbook's default (unstyled) look before writing values.
    $ws.Rows.Item($r).ClearFormats()
}

# New row data: row number, contenedor, fecha (date serial), comuna, empresa, servicios
$newData = @(
    @(2,  "MSKU9167801", 45174.6875,         "San Antonio", "sti", 82589),
    @(3,  "MSKU8950767", 45174.69791666666,  "San Antonio", "sti", 82566),
    @(4,  "MRKU2590580", 45174.70833333334,  "San Antonio", "sti", 82571),
    @(5,  "HASU4761072", 45174.73958333334,  "San Antonio", "sti", 82590),
    @(6,  "MSKU1728038", 45174.83333333334,  "San Antonio", "sti", 82152),
    @(9,  "AMCU9389196", 45174.75347222222,  "San Antonio", "sti", 83961),
    @(40, "MEDU1070073", 45174.83888888889,  "San Antonio", "sti", 83389)
)

foreach ($row in $newData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $row[2]
    $cellB.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
